# IndividualWorkSheet.xlsx - "Add files via upload" update
# Fills in the git-commit / follow-up progress data for 유병주's rows 7-8
# (Member.java / Menu.java git commits, and the Table.java GUI bug note),
# and leaves the workbook with that sheet active/selected, matching the
# state it was re-uploaded in.

$wb = $excel.ActiveWorkbook

# 유병주 is the 3rd sheet (rId3 -> xl/worksheets/sheet3.xml)
$ws = $wb.Worksheets.Item(3)

# Row 7: "Java와 Mysql 연동하기(2)" task finished 2019-05-22 with a note
# about the Member.java commit.
$ws.Range("D7").Value = "2019-05-22"
$ws.Range("E7").Value = "수정한 Member.java git commit"

# Row 8: "Java와 Mysql 연동하기(4)" task finished 2019-05-22, noting the
# Menu.java commit plus a problem discovered in the Table.java GUI.
$ws.Range("D8").Value = "2019-05-22"
$ws.Range("E8").Value = "수정한 Menu.java git commit"
$ws.Range("F8").Value = "GUI Table.java에서 주문내역에 Jtable에 추가되지 않는 문제점이 발견 -> Table 주문내역을 데이터베이스에 저장하는데 시간이 필요"

# Switch the active/selected tab from 박서린 to 유병주, with the last
# selection sitting at L8.
$ws.Activate()
[void]$ws.Range("L8").Select()
